$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (2nd paragraph, right after the title)
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new paragraph before the last paragraph (the "Prompt:" paragraph)
#    that holds the bold title text, then update the last paragraph's text.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range
$insertRange.Collapse(1)  # wdCollapseStart
$insertRange.InsertParagraphBefore()

# Now re-fetch: the "prompt" paragraph moved down by one; the newly created
# paragraph is the second-to-last paragraph.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara.Range.Text = "Play Aloha Fruit Bonanza Free: Slot Review & RTP"
$newPara.Range.Font.Bold = 1

# 3. Replace text of the last ("Prompt:") paragraph and make it italic
$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptPara.Range.Text = "Review of Aloha Fruit Bonanza slot: gameplay, payouts, free spins. Play for free and win up to 7,500x the stake. RTP of 97.01%."
$promptPara.Range.Font.Italic = 1
